$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking price strings are not
# auto-converted to floating point numbers (the source cells are plain text).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.262.55'
$ws.Range("E2").Value = '  +0.96%  '
$ws.Range("D3").Value = '1.676.88'
$ws.Range("E3").Value = '  +0.62%  '
$ws.Range("E4").Value = '  +0.27%  '
$ws.Range("D5").Value = '217.25'
$ws.Range("E5").Value = '  +0.34%  '
$ws.Range("D6").Value = '0.5294'
$ws.Range("D7").Value = '1.007'
$ws.Range("E7").Value = '  +0.26%  '
$ws.Range("D8").Value = '0.2684'
$ws.Range("E8").Value = '  +1.72%  '
$ws.Range("D9").Value = '0.06472'
$ws.Range("E9").Value = '  +1.27%  '
$ws.Range("E10").Value = '  -0.13%  '
$ws.Range("D11").Value = '0.07508'
$ws.Range("E11").Value = '  +1.25%  '
$ws.Range("D12").Value = '1.672.20'
$ws.Range("E12").Value = '  +0.23%  '
$ws.Range("D13").Value = '4.511'
$ws.Range("E13").Value = '  +0.28%  '
$ws.Range("D14").Value = '0.5764'
$ws.Range("E14").Value = '  -1.04%  '
$ws.Range("D15").Value = '0.000008492'
$ws.Range("E15").Value = '  +0.14%  '
$ws.Range("D16").Value = '64.59'
$ws.Range("E16").Value = '  +0.67%  '
$ws.Range("D17").Value = '26.298.13'
$ws.Range("E17").Value = '  +0.93%  '
$ws.Range("D18").Value = '4.907'
$ws.Range("E18").Value = '  -0.38%  '
$ws.Range("E19").Value = '  +0.24%  '
$ws.Range("E20").Value = '  +1.20%  '
$ws.Range("D21").Value = '190.08'
$ws.Range("E21").Value = '  +0.25%  '
$ws.Range("E22").Value = '  -0.43%  '
$ws.Range("E23").Value = '  +0.25%  '
$ws.Range("D24").Value = '144.92'
$ws.Range("E24").Value = '  -0.06%  '
$ws.Range("D25").Value = '7.795'
$ws.Range("E25").Value = '  +2.53%  '
$ws.Range("D26").Value = '0.1265'
$ws.Range("E26").Value = '  +6.03%  '
$ws.Range("D27").Value = '15.74'
$ws.Range("E27").Value = '  +0.67%  '
$ws.Range("D28").Value = '0.06505'
$ws.Range("E28").Value = '  -2.09%  '
$ws.Range("D29").Value = '1.363'
$ws.Range("E29").Value = '  +3.90%  '
$ws.Range("E30").Value = '  +0.34%  '
$ws.Range("D31").Value = '3.584'
$ws.Range("E31").Value = '  +1.84%  '
$ws.Range("E32").Value = '  +2.29%  '
$ws.Range("D33").Value = '1.652'
$ws.Range("E33").Value = '  +1.14%  '
$ws.Range("D34").Value = '1.030'
$ws.Range("E34").Value = '  +1.37%  '
$ws.Range("D35").Value = '0.6189'
$ws.Range("E35").Value = '  +1.81%  '
$ws.Range("D36").Value = '2.403'
$ws.Range("E36").Value = '  +1.52%  '
$ws.Range("D37").Value = '2.728'
$ws.Range("E37").Value = '  +0.68%  '
$ws.Range("E38").Value = '  +1.18%  '
$ws.Range("D39").Value = '1.113.89'
$ws.Range("E39").Value = '  +3.59%  '
$ws.Range("D40").Value = '0.01621'
$ws.Range("E40").Value = '  +0.95%  '
$ws.Range("D41").Value = '0.8726'
$ws.Range("E41").Value = '  +1.81%  '
$ws.Range("D42").Value = '1.015'
$ws.Range("E42").Value = '  +0.66%  '
$ws.Range("D43").Value = '100.44'
$ws.Range("E43").Value = '  +0.00%  '
$ws.Range("D44").Value = '1.827.97'
$ws.Range("E44").Value = '  +0.77%  '
$ws.Range("D45").Value = '56.90'
$ws.Range("E45").Value = '  +1.18%  '
$ws.Range("D46").Value = '8.176'
$ws.Range("E46").Value = '  +1.93%  '
$ws.Range("E47").Value = '  +0.06%  '
$ws.Range("E48").Value = '  +1.02%  '
$ws.Range("D49").Value = '0.4290'
$ws.Range("E49").Value = '  +0.05%  '
$ws.Range("D50").Value = '6.073'
$ws.Range("E50").Value = '  +1.98%  '
$ws.Range("D51").Value = '0.3373'
$ws.Range("E51").Value = '  +1.93%  '

# Restore the default (unstyled) cell style on column D now that the text
# values are written, so no stray style index lingers on the cells.
$ws.Range("D2:D51").Style = "Normal"
